# Applies: "Do 400 ao 500 classificados"
# Fills in the classification labels (column B) for rows 400-500 on the
# "Treinamento" sheet, widens column A, and restores the view state
# (scroll position / zoom / selection) that Excel saved after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

$values = @(1,1,0,1,1,0,0,0,1,0,1,1,1,0,0,1,0,0,1,1,0,0,0,0,0,0,0,1,0,0,0,1,0,0,0,0,1,1,1,1,0,0,0,0,0,1,1,0,1,1,1,1,0,0,1,1,0,1,1,1,0,1,0,1,0,1,1,1,1,1,1,0,1,1,1,1,0,0,1,1,1,0,0,1,1,0,1,1,1,1,1,1,1,1,0,1,1,1,0,0,1)

$startRow = 400
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Column A is wide enough to read full tweet text.
$ws.Columns.Item(1).ColumnWidth = 255.43

# Restore view: scrolled near the bottom, zoomed to 70%, with B501 selected.
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("A378").Select()
$ws.Range("B501").Select()
